$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text values (shared-string content changes)
$ws.Range("I1").Value = "time"
$ws.Range("A1").Value = "Speed"
$ws.Range("B1").Value = "Draught"

# Select A20:I21 (active cell A20), matching the recorded selection
$ws.Range("A20:I21").Select()

# Best-effort: scroll the view so A4 becomes the top-left visible cell
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
